$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated accuracy values for B2:B114 (MO4 Froze Token Embeddings Layers re-run)
$values = @(0.984375,0.96875,0.96875,0.96875,0.96875,0.953125,0.953125,0.953125,0.90625,0.953125,0.921875,0.953125,1,0.921875,0.953125,0.90625,0.9375,0.90625,0.921875,0.875,0.90625,0.890625,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.875,0.890625,0.890625,0.890625,0.890625,0.890625,0.890625,0.890625,0.890625,0.890625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.90625,0.921875,0.9375,0.953125,0.921875,0.921875,0.921875,0.984375,0.921875,0.9375,0.953125,1)

$rows = $values.Length
$arr = New-Object 'object[,]' $rows,1
for ($i = 0; $i -lt $rows; $i++) {
    $arr[$i,0] = $values[$i]
}
$ws.Range("B2:B114").Value = $arr

# Refresh the repr strings in column A for the re-run DisplayOutputs object (new memory address)
$newRepr = "<__main__.DisplayOutputs object at 0x7f3ebc643100>"
$ws.Range("A102:A114").Value = $newRepr

# Update the sheet view / selection to match the saved state
$ws.Range("A2:B114").Select()
